$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "24.457.87"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.667.99"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3964"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3918"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.05"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.407"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08606"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.367"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("E15").Value = "  +4.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.929"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.03%  "
$ws.Range("D17").Value = "1.668.20"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06983"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.027"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "24.450.66"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.426"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.032"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.77%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.494"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "143.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.151"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.549"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.17%  "
$ws.Range("D33").Value = "1.849.68"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.065"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08288"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03029"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.867"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2773"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09263"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.445"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7153"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.546"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.147"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08456"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.291"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.64%  "
